# Generate Report for Handoff
#
# This reflects a new handoff run: the "10a56bb5-9f2d-4704-a2b9-9afa2331fb10"
# file (rows 7,8,9,10,12,14 in each per-language table -- row 11/13 entries
# were not part of this handoff) now carries priority "ht" and refreshed
# handoff timestamps, and the Overview sheet's "Latest HO Xliff Generate
# Date" column picks up the same refreshed timestamp.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 12, 14)

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-09-04 18:25:11"
}

# zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-09-04 18:25:02"
}

# de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
# (the de-de table's "Latest Handoff Datetime" happened to share the exact
# same placeholder timestamp string as the Overview sheet, so it refreshes
# alongside the Overview column in this handoff run as well)
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-09-04 18:25:11"
}
